$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Seer" shifts to D, "Healer" shifts to E).
$ws.Columns("C:C").Insert()

# Rows 2 and 3 ("Perception boost" / "See status") are not archetype-specific columns
# that moved -- they stay put in column C (only their shared-string index changes).
$ws.Range("C2").Value = $ws.Range("D2").Value2
$ws.Range("D2").ClearContents()
$ws.Range("C3").Value = $ws.Range("D3").Value2
$ws.Range("D3").ClearContents()

# New header for the inserted "Empath" column.
$ws.Range("C1").Value = "Empath"

# The inserted column keeps the wide "description" column width (same as column D).
$ws.Columns("C").ColumnWidth = 40.5

# Update the active selection to match the edited workbook's last cursor position.
$null = $ws.Range("C3").Select()
